$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 9).Value = 0.977669497583861
$ws.Cells.Item(2, 10).Value = 0.977669497583861
$ws.Cells.Item(2, 13).Value = 1.363113
$ws.Cells.Item(2, 14).Value = 4.089339
$ws.Cells.Item(2, 15).Value = 0.06061833851125786
$ws.Cells.Item(2, 16).Value = 0.06061833851125786
$ws.Cells.Item(2, 17).Value = 10.243514756835
$ws.Cells.Item(2, 18).Value = 92.19163281151499
$ws.Cells.Item(2, 19).Value = 0.05926470055666989
$ws.Cells.Item(2, 20).Value = 0.05926470055666988

$ws.Cells.Item(3, 9).Value = 0.977669497583861
$ws.Cells.Item(3, 10).Value = 0.977669497583861
$ws.Cells.Item(3, 15).Value = 0.1775969932713293
$ws.Cells.Item(3, 16).Value = 0.1775969932713292
$ws.Cells.Item(3, 19).Value = 0.1736311631839848
$ws.Cells.Item(3, 20).Value = 0.1736311631839848

$ws.Cells.Item(4, 9).Value = 0.977669497583861
$ws.Cells.Item(4, 10).Value = 0.977669497583861
$ws.Cells.Item(4, 13).Value = 3.206217333333333
$ws.Cells.Item(4, 14).Value = 9.618651999999999
$ws.Cells.Item(4, 15).Value = 0.1425821393036839
$ws.Cells.Item(4, 16).Value = 0.1425821393036839
$ws.Cells.Item(4, 17).Value = 24.09406598544666
$ws.Cells.Item(4, 18).Value = 216.84659386902
$ws.Cells.Item(4, 19).Value = 0.1393982084974647
$ws.Cells.Item(4, 20).Value = 0.1393982084974647

$ws.Cells.Item(5, 9).Value = 0.977669497583861
$ws.Cells.Item(5, 10).Value = 0.977669497583861
$ws.Cells.Item(5, 13).Value = 2.765104
$ws.Cells.Item(5, 14).Value = 8.295312000000001
$ws.Cells.Item(5, 15).Value = 0.1229656017445606
$ws.Cells.Item(5, 16).Value = 0.1229656017445605
$ws.Cells.Item(5, 17).Value = 20.77918971368
$ws.Cells.Item(5, 18).Value = 187.01270742312
$ws.Cells.Item(5, 19).Value = 0.1202197180777017
$ws.Cells.Item(5, 20).Value = 0.1202197180777017

$ws.Cells.Item(6, 9).Value = 0.977669497583861
$ws.Cells.Item(6, 10).Value = 0.977669497583861
$ws.Cells.Item(6, 13).Value = 4.898567333333333
$ws.Cells.Item(6, 14).Value = 14.695702
$ws.Cells.Item(6, 15).Value = 0.2178418170996753
$ws.Cells.Item(6, 16).Value = 0.2178418170996753
$ws.Cells.Item(6, 17).Value = 36.81172930369667
$ws.Cells.Item(6, 18).Value = 331.30556373327
$ws.Cells.Item(6, 19).Value = 0.2129772998765949
$ws.Cells.Item(6, 20).Value = 0.2129772998765949

$ws.Cells.Item(7, 9).Value = 0.977669497583861
$ws.Cells.Item(7, 10).Value = 0.977669497583861
$ws.Cells.Item(7, 13).Value = 6.260217666666667
$ws.Cells.Item(7, 14).Value = 18.780653
$ws.Cells.Item(7, 15).Value = 0.278395110069493
$ws.Cells.Item(7, 16).Value = 0.278395110069493
$ws.Cells.Item(7, 17).Value = 47.04425242037833
$ws.Cells.Item(7, 18).Value = 423.398271783405
$ws.Cells.Item(7, 19).Value = 0.2721784073914449
$ws.Cells.Item(7, 20).Value = 0.2721784073914449

$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.171642
$ws.Cells.Item(8, 8).Value = 0.514926
$ws.Cells.Item(8, 9).Value = 0.02233050241613897
$ws.Cells.Item(8, 10).Value = 0.02233050241613898
$ws.Cells.Item(8, 13).Value = 1.363113
$ws.Cells.Item(8, 14).Value = 4.089339
$ws.Cells.Item(8, 15).Value = 0.06061833851125786
$ws.Cells.Item(8, 16).Value = 0.06061833851125786
$ws.Cells.Item(8, 17).Value = 0.233967441546
$ws.Cells.Item(8, 18).Value = 2.105706973914
$ws.Cells.Item(8, 19).Value = 0.001353637954587974
$ws.Cells.Item(8, 20).Value = 0.001353637954587974

$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.171642
$ws.Cells.Item(9, 8).Value = 0.514926
$ws.Cells.Item(9, 9).Value = 0.02233050241613897
$ws.Cells.Item(9, 10).Value = 0.02233050241613898
$ws.Cells.Item(9, 15).Value = 0.1775969932713293
$ws.Cells.Item(9, 16).Value = 0.1775969932713292
$ws.Cells.Item(9, 17).Value = 0.685467717566
$ws.Cells.Item(9, 18).Value = 6.169209458094
$ws.Cells.Item(9, 19).Value = 0.003965830087344436
$ws.Cells.Item(9, 20).Value = 0.003965830087344436

$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.171642
$ws.Cells.Item(10, 8).Value = 0.514926
$ws.Cells.Item(10, 9).Value = 0.02233050241613897
$ws.Cells.Item(10, 10).Value = 0.02233050241613898
$ws.Cells.Item(10, 13).Value = 3.206217333333333
$ws.Cells.Item(10, 14).Value = 9.618651999999999
$ws.Cells.Item(10, 15).Value = 0.1425821393036839
$ws.Cells.Item(10, 16).Value = 0.1425821393036839
$ws.Cells.Item(10, 17).Value = 0.5503215555279999
$ws.Cells.Item(10, 18).Value = 4.952893999752
$ws.Cells.Item(10, 19).Value = 0.003183930806219177
$ws.Cells.Item(10, 20).Value = 0.003183930806219177

$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.171642
$ws.Cells.Item(11, 8).Value = 0.514926
$ws.Cells.Item(11, 9).Value = 0.02233050241613897
$ws.Cells.Item(11, 10).Value = 0.02233050241613898
$ws.Cells.Item(11, 13).Value = 2.765104
$ws.Cells.Item(11, 14).Value = 8.295312000000001
$ws.Cells.Item(11, 15).Value = 0.1229656017445606
$ws.Cells.Item(11, 16).Value = 0.1229656017445605
$ws.Cells.Item(11, 17).Value = 0.474607980768
$ws.Cells.Item(11, 18).Value = 4.271471826912
$ws.Cells.Item(11, 19).Value = 0.002745883666858893
$ws.Cells.Item(11, 20).Value = 0.002745883666858893

$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.171642
$ws.Cells.Item(12, 8).Value = 0.514926
$ws.Cells.Item(12, 9).Value = 0.02233050241613897
$ws.Cells.Item(12, 10).Value = 0.02233050241613898
$ws.Cells.Item(12, 13).Value = 4.898567333333333
$ws.Cells.Item(12, 14).Value = 14.695702
$ws.Cells.Item(12, 15).Value = 0.2178418170996753
$ws.Cells.Item(12, 16).Value = 0.2178418170996753
$ws.Cells.Item(12, 17).Value = 0.8407998942279999
$ws.Cells.Item(12, 18).Value = 7.567199048052
$ws.Cells.Item(12, 19).Value = 0.004864517223080404
$ws.Cells.Item(12, 20).Value = 0.004864517223080404

$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.171642
$ws.Cells.Item(13, 8).Value = 0.514926
$ws.Cells.Item(13, 9).Value = 0.02233050241613897
$ws.Cells.Item(13, 10).Value = 0.02233050241613898
$ws.Cells.Item(13, 13).Value = 6.260217666666667
$ws.Cells.Item(13, 14).Value = 18.780653
$ws.Cells.Item(13, 15).Value = 0.278395110069493
$ws.Cells.Item(13, 16).Value = 0.278395110069493
$ws.Cells.Item(13, 17).Value = 1.074516280742
$ws.Cells.Item(13, 18).Value = 9.670646526678
$ws.Cells.Item(13, 19).Value = 0.006216702678048089
$ws.Cells.Item(13, 20).Value = 0.00621670267804809
